$d = $word.ActiveDocument

# 1. Apply single line spacing (w:line="240" w:lineRule="auto") to every
#    paragraph in the main document story. This covers all the body
#    paragraphs (including the empty "<w:p/>" spacer paragraphs between
#    the PERTAMA/KEDUA/KETIGA/KEEMPAT diktum items).
$d.Content.ParagraphFormat.LineSpacingRule = 0

# 2. The table-cell paragraphs at the bottom of the document (the
#    signature block) are not reachable through the bulk
#    ParagraphFormat assignment above, so set them explicitly by
#    reinserting each cell's first paragraph with the same content plus
#    the new <w:spacing> setting, via Range.InsertXML (which operates
#    precisely on the paragraph's own range).
$t = $d.Tables.Item(1)

function Set-CellParaSpacing($row, $col, $innerXml) {
    $cell = $t.Cell($row, $col)
    $p = $cell.Range.Paragraphs.Item(1)
    $r = $p.Range
    $xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:line="240" w:lineRule="auto"/></w:pPr>' + $innerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($xmlFrag) | Out-Null
}

Set-CellParaSpacing 1 1 ""
Set-CellParaSpacing 1 2 "<w:r><w:t>Ditetapkan di : {{satker_kota}}</w:t></w:r>"
Set-CellParaSpacing 2 1 ""
Set-CellParaSpacing 2 2 "<w:r><w:t>Pada tanggal  : {{tanggal_sk}}</w:t></w:r>"
Set-CellParaSpacing 3 1 ""
Set-CellParaSpacing 3 2 "<w:r/>"
Set-CellParaSpacing 4 1 ""
Set-CellParaSpacing 4 2 "<w:r><w:t>KUASA PENGGUNA ANGGARAN,</w:t></w:r>"
Set-CellParaSpacing 5 1 ""
Set-CellParaSpacing 5 2 "<w:r><w:br/><w:br/><w:br/><w:br/><w:t>{{kpa_nama}}</w:t><w:br/><w:t>NIP. {{kpa_nip}}</w:t></w:r>"

Write-Output "Line spacing set to single for all paragraphs."
